$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 9.881153648301549
$ws.Range("C2").Value = 7.915017854788083
$ws.Range("D2").Value = 3.733362436769495
$ws.Range("F2").Value = 19.9866209883833
$ws.Range("G2").Value = 3.5936690653342
$ws.Range("I2").Value = 17.0461601231509
$ws.Range("M2").Value = 19.96700386071196
$ws.Range("N2").Value = 16.26924410296094
$ws.Range("O2").Value = 17.6917833093795

$ws.Range("B3").Value = 9.373212588027036
$ws.Range("C3").Value = 7.584907697067863
$ws.Range("D3").Value = 3.675316739641407
$ws.Range("F3").Value = 19.90738053844318
$ws.Range("G3").Value = 3.595554913489552
$ws.Range("I3").Value = 17.11957829895055
$ws.Range("M3").Value = 19.35892791304355
$ws.Range("N3").Value = 16.33121745991418
$ws.Range("O3").Value = 17.69843233823755

$ws.Range("B4").Value = 9.047258988756798
$ws.Range("C4").Value = 7.373159530333133
$ws.Range("D4").Value = 3.638715676524649
$ws.Range("F4").Value = 19.86523609738898
$ws.Range("G4").Value = 3.596774418860985
$ws.Range("I4").Value = 17.168856490374
$ws.Range("M4").Value = 18.98378708692832
$ws.Range("N4").Value = 16.37105263929876
$ws.Range("O4").Value = 17.7075900206965

$ws.Range("B5").Value = 8.911034962097082
$ws.Range("C5").Value = 7.284668387007853
$ws.Range("D5").Value = 3.623568412166338
$ws.Range("F5").Value = 19.84971250780678
$ws.Range("G5").Value = 3.597286914739892
$ws.Range("I5").Value = 17.18999015601645
$ws.Range("M5").Value = 18.83072535755881
$ws.Range("N5").Value = 16.38773595312628
$ws.Range("O5").Value = 17.71259592566362

$ws.Range("B6").Value = 8.888214722258935
$ws.Range("C6").Value = 7.269843938531759
$ws.Range("D6").Value = 3.62103949229448
$ws.Range("F6").Value = 19.84723487211057
$ws.Range("G6").Value = 3.597372954218551
$ws.Range("I6").Value = 17.19356284477089
$ws.Range("M6").Value = 18.80530598917477
$ws.Range("N6").Value = 16.39053344379492
$ws.Range("O6").Value = 17.71350402550485

$ws.Range("B7").Value = 9.045435361882323
$ws.Range("C7").Value = 7.371974915434203
$ws.Range("D7").Value = 3.638512321850156
$ws.Range("F7").Value = 19.86502004134503
$ws.Range("G7").Value = 3.596781267581049
$ws.Range("I7").Value = 17.16913724993816
$ws.Range("M7").Value = 18.98172325587807
$ws.Range("N7").Value = 16.37127581140619
$ws.Range("O7").Value = 17.70765237684732

$ws.Range("B8").Value = 9.709015808404226
$ws.Range("C8").Value = 7.803118363362241
$ws.Range("D8").Value = 3.713553311639583
$ws.Range("F8").Value = 19.95795684552901
$ws.Range("G8").Value = 3.594306554157835
$ws.Range("I8").Value = 17.07060108080558
$ws.Range("M8").Value = 19.75785961051755
$ws.Range("N8").Value = 16.29024344126664
$ws.Range("O8").Value = 17.69302151879839

$ws.Range("B9").Value = 10.89346847218646
$ws.Range("C9").Value = 8.57389669656397
$ws.Range("D9").Value = 3.852684650337248
$ws.Range("F9").Value = 20.19113234598985
$ws.Range("G9").Value = 3.589939985651769
$ws.Range("I9").Value = 16.91084597122427
$ws.Range("M9").Value = 21.25511578786506
$ws.Range("N9").Value = 16.14540915751734
$ws.Range("O9").Value = 17.70465840762173

$ws.Range("B10").Value = 11.74066449909678
$ws.Range("C10").Value = 9.091534370923538
$ws.Range("D10").Value = 3.949505157808151
$ws.Range("F10").Value = 20.3923572729509
$ws.Range("G10").Value = 3.587025090936742
$ws.Range("I10").Value = 16.81409228994817
$ws.Range("M10").Value = 22.32700807729825
$ws.Range("N10").Value = 16.04746446709638
$ws.Range("O10").Value = 17.73782123256647

$ws.Range("B11").Value = 12.11539766440761
$ws.Range("C11").Value = 9.315956523545706
$ws.Range("D11").Value = 3.992279503861704
$ws.Range("F11").Value = 20.4901046883656
$ws.Range("G11").Value = 3.585762005210024
$ws.Range("I11").Value = 16.77459850789732
$ws.Range("M11").Value = 22.80612216390144
$ws.Range("N11").Value = 16.00472094959486
$ws.Range("O11").Value = 17.75823992887898

$ws.Range("B12").Value = 12.25403268317952
$ws.Range("C12").Value = 9.399316671965243
$ws.Range("D12").Value = 4.008286649080699
$ws.Range("F12").Value = 20.52798441987555
$ws.Range("G12").Value = 3.585292701560291
$ws.Range("I12").Value = 16.76029688531065
$ws.Range("M12").Value = 22.98615255808928
$ws.Range("N12").Value = 15.98879383049911
$ws.Range("O12").Value = 17.76673617984455

$ws.Range("B13").Value = 12.22432050035114
$ws.Range("C13").Value = 9.381436268597982
$ws.Range("D13").Value = 4.00484782202207
$ws.Range("F13").Value = 20.5197883400636
$ws.Range("G13").Value = 3.585393374996717
$ws.Range("I13").Value = 16.76334785614922
$ws.Range("M13").Value = 22.94744469030635
$ws.Range("N13").Value = 15.99221253094522
$ws.Range("O13").Value = 17.7648724281912

$ws.Range("B14").Value = 12.12686874973778
$ws.Range("C14").Value = 9.32284732897576
$ws.Range("D14").Value = 3.993600284803717
$ws.Range("F14").Value = 20.49320390948448
$ws.Range("G14").Value = 3.585723215224601
$ws.Range("I14").Value = 16.77340877946681
$ws.Range("M14").Value = 22.82096238177948
$ws.Range("N14").Value = 16.00340543609898
$ws.Range("O14").Value = 17.75892362502162

$ws.Range("B15").Value = 12.06675112002042
$ws.Range("C15").Value = 9.286747548226531
$ws.Range("D15").Value = 3.986685796193104
$ws.Range("F15").Value = 20.4770319524987
$ws.Range("G15").Value = 3.585926422483023
$ws.Range("I15").Value = 16.77965664511138
$ws.Range("M15").Value = 22.74330105106771
$ws.Range("N15").Value = 16.01029508452039
$ws.Range("O15").Value = 17.75537922772255

$ws.Range("B16").Value = 11.7157174587403
$ws.Range("C16").Value = 9.076642376169914
$ws.Range("D16").Value = 3.946683510411024
$ws.Range("F16").Value = 20.38609184693555
$ws.Range("G16").Value = 3.587108898490629
$ws.Range("I16").Value = 16.8167645847839
$ws.Range("M16").Value = 22.2955102647263
$ws.Range("N16").Value = 16.05029417831474
$ws.Range("O16").Value = 17.73659392110826

$ws.Range("B17").Value = 11.49454887971348
$ws.Range("C17").Value = 8.9448923818547
$ws.Range("D17").Value = 3.921812366087783
$ws.Range("F17").Value = 20.33187407170877
$ws.Range("G17").Value = 3.587850389229877
$ws.Range("I17").Value = 16.84068952678109
$ws.Range("M17").Value = 22.01849686595246
$ws.Range("N17").Value = 16.07529525639644
$ws.Range("O17").Value = 17.7264339594954

$ws.Range("B18").Value = 11.36520001338535
$ws.Range("C18").Value = 8.868074478522313
$ws.Range("D18").Value = 3.907388072584453
$ws.Range("F18").Value = 20.3012756805662
$ws.Range("G18").Value = 3.588282799727558
$ws.Range("I18").Value = 16.85487573242782
$ws.Range("M18").Value = 21.85837598137297
$ws.Range("N18").Value = 16.08984586832333
$ws.Range("O18").Value = 17.72109238469666

$ws.Range("B19").Value = 11.32103750279486
$ws.Range("C19").Value = 8.841887977456309
$ws.Range("D19").Value = 3.902484054636637
$ws.Range("F19").Value = 20.29101709893548
$ws.Range("G19").Value = 3.588430225506125
$ws.Range("I19").Value = 16.85975185029057
$ws.Range("M19").Value = 21.80403182242192
$ws.Range("N19").Value = 16.09480181543385
$ws.Range("O19").Value = 17.71937013406432

$ws.Range("B20").Value = 11.51831400339111
$ws.Range("C20").Value = 8.959025169888786
$ws.Range("D20").Value = 3.924472324230871
$ws.Range("F20").Value = 20.3375851694929
$ws.Range("G20").Value = 3.587770843443018
$ws.Range("I20").Value = 16.83809864002426
$ws.Range("M20").Value = 22.04806853531593
$ws.Range("N20").Value = 16.07261619893245
$ws.Range("O20").Value = 17.72746354912363

$ws.Range("B21").Value = 12.15558138806182
$ws.Range("C21").Value = 9.340100630236739
$ws.Range("D21").Value = 3.996909192583086
$ws.Range("F21").Value = 20.50098915997626
$ws.Range("G21").Value = 3.585626089287357
$ws.Range("I21").Value = 16.77043586757052
$ws.Range("M21").Value = 22.85815262043285
$ws.Range("N21").Value = 16.00011079448101
$ws.Range("O21").Value = 17.76065022056259

$ws.Range("B22").Value = 12.5530236399366
$ws.Range("C22").Value = 9.579679587763158
$ws.Range("D22").Value = 4.043136686656112
$ws.Range("F22").Value = 20.61281020212141
$ws.Range("G22").Value = 3.584276802585919
$ws.Range("I22").Value = 16.73002702485237
$ws.Range("M22").Value = 23.37936069021312
$ws.Range("N22").Value = 15.95423276690081
$ws.Range("O22").Value = 17.78679167111735

$ws.Range("B23").Value = 12.3426439138011
$ws.Range("C23").Value = 9.452688686117439
$ws.Range("D23").Value = 4.018568673545129
$ws.Range("F23").Value = 20.55267894097267
$ws.Range("G23").Value = 3.584992160197952
$ws.Range("I23").Value = 16.75124386281492
$ws.Range("M23").Value = 23.10198931170088
$ws.Range("N23").Value = 15.97858125161748
$ws.Range("O23").Value = 17.77243326082542

$ws.Range("B24").Value = 11.50757663121828
$ws.Range("C24").Value = 8.952639076164932
$ws.Range("D24").Value = 3.923270147577709
$ws.Range("F24").Value = 20.33500139802924
$ws.Range("G24").Value = 3.587806787035164
$ws.Range("I24").Value = 16.83926863638688
$ws.Range("M24").Value = 22.03470185542003
$ws.Range("N24").Value = 16.0738268488766
$ws.Range("O24").Value = 17.72699651555741

$ws.Range("B25").Value = 10.58618447983003
$ws.Range("C25").Value = 8.373731384842529
$ws.Range("D25").Value = 3.815958925561931
$ws.Range("F25").Value = 20.12270804555521
$ws.Range("G25").Value = 3.591069530541306
$ws.Range("I25").Value = 17.0461601231509
$ws.Range("M25").Value = 20.85413298547299
$ws.Range("N25").Value = 16.18309602920111
$ws.Range("O25").Value = 17.69718648789772
